$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (A:F narrow, G wide "Notes" column)
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").EntireColumn.ColumnWidth = 11.65
$ws.Range("G1").EntireColumn.ColumnWidth = 43.65

# ---------------------------------------------------------------------------
# Header row (row 1): Event Path | Timeline | 3D | Parameter | Para. Type | Para. Range | Notes
# Bold black font, full thin black box border around every cell.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Event Path"
$ws.Range("B1").Value = "Timeline"
$ws.Range("C1").Value = "3D"
$ws.Range("D1").Value = "Parameter"
$ws.Range("E1").Value = "Para. Type"
$ws.Range("F1").Value = "Para. Range"
$ws.Range("G1").Value = "Notes"

$hdr = $ws.Range("A1")
$hdr.Font.Bold = $true
$hdr.Font.Color = 0
$hdr.Borders.Color = 0
$hdr.Borders.Weight = 2
$hdr.Borders.LineStyle = 1
$hdr.Copy() | Out-Null
$ws.Range("B1:G1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Rows 2-4: plain body rows, full thin black box border, regular black font.
# Column C = "N", Column A = ":/" ; other columns blank but still styled.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "N"
$ws.Range("A2").Value = ":/"
$ws.Range("C3").Value = "N"
$ws.Range("A3").Value = ":/"
$ws.Range("C4").Value = "N"
$ws.Range("A4").Value = ":/"

$body1 = $ws.Range("A2")
$body1.Font.Bold = $false
$body1.Font.Color = 0
$body1.Borders.Color = 0
$body1.Borders.Weight = 2
$body1.Borders.LineStyle = 1
$body1.Copy() | Out-Null
$ws.Range("B2:G4").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Rows 5-8: same text pattern, but bordered with left/right/top only (no
# bottom), and the "Notes" column (G) additionally has wrap text enabled.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "N"
$ws.Range("A5").Value = ":/"
$ws.Range("C6").Value = "N"
$ws.Range("A6").Value = ":/"
$ws.Range("C7").Value = "N"
$ws.Range("A7").Value = ":/"
$ws.Range("C8").Value = "N"
$ws.Range("A8").Value = ":/"

$body2 = $ws.Range("A5")
$body2.Font.Bold = $false
$body2.Font.Color = 0
$body2.Borders.Color = 0
$body2.Borders.Weight = 2
$body2.Borders.Item(7).LineStyle = 1
$body2.Borders.Item(8).LineStyle = 1
$body2.Borders.Item(10).LineStyle = 1
$body2.Copy() | Out-Null
$ws.Range("B5:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:A8").PasteSpecial(-4122) | Out-Null

$notes = $ws.Range("G5")
$notes.Font.Bold = $false
$notes.Font.Color = 0
$notes.Borders.Color = 0
$notes.Borders.Weight = 2
$notes.Borders.Item(7).LineStyle = 1
$notes.Borders.Item(8).LineStyle = 1
$notes.Borders.Item(10).LineStyle = 1
$notes.WrapText = $true
$notes.Copy() | Out-Null
$ws.Range("G6:G8").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("D16").Select() | Out-Null
